$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '59.654.13'
$c.ClearFormats()
$ws.Range("E2").Value = '  +3.96%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.018.47'
$c.ClearFormats()
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("E4").Value = '  +0.08%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '564.69'
$c.ClearFormats()
$ws.Range("E5").Value = '  +3.29%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '139.72'
$c.ClearFormats()
$ws.Range("E6").Value = '  +8.25%  '
$ws.Range("E7").Value = '  +0.01%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.520'
$c.ClearFormats()
$ws.Range("E8").Value = '  +2.24%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '3.008.78'
$c.ClearFormats()
$ws.Range("E9").Value = '  +3.16%  '
$ws.Range("E10").Value = '  +6.29%  '
$ws.Range("E11").Value = '  +10.97%  '
$ws.Range("E12").Value = '  +3.43%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000232'
$c.ClearFormats()
$ws.Range("E13").Value = '  +5.89%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '33.98'
$c.ClearFormats()
$ws.Range("E14").Value = '  +3.89%  '
$ws.Range("E15").Value = '  +1.92%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.516.96'
$c.ClearFormats()
$ws.Range("E16").Value = '  +3.44%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '7.26'
$c.ClearFormats()
$ws.Range("E17").Value = '  +6.61%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.019.31'
$c.ClearFormats()
$ws.Range("E18").Value = '  +3.51%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '59.648.16'
$c.ClearFormats()
$ws.Range("E19").Value = '  +3.91%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '435.42'
$c.ClearFormats()
$ws.Range("E20").Value = '  +4.66%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '13.67'
$c.ClearFormats()
$ws.Range("E21").Value = '  +4.42%  '
$ws.Range("E22").Value = '  +6.31%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '7.13'
$c.ClearFormats()
$ws.Range("E23").Value = '  +2.83%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '13.26'
$c.ClearFormats()
$ws.Range("E24").Value = '  +2.26%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '80.81'
$c.ClearFormats()
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("E26").Value = '  +0.13%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.ClearFormats()
$ws.Range("E27").Value = '  +13.99%  '
$ws.Range("E28").Value = '  +0.13%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.55'
$c.ClearFormats()
$ws.Range("E29").Value = '  +3.78%  '
$ws.Range("E30").Value = '  +5.90%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '26.02'
$c.ClearFormats()
$ws.Range("E31").Value = '  +3.55%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.27'
$c.ClearFormats()
$ws.Range("E32").Value = '  +5.83%  '
$ws.Range("E33").Value = '  +5.73%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0₃0783'
$c.ClearFormats()
$ws.Range("E34").Value = '  +16.38%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E35").Value = '  +7.71%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '5.93'
$c.ClearFormats()
$ws.Range("E36").Value = '  +5.55%  '
$ws.Range("E37").Value = '  +3.06%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '49.17'
$c.ClearFormats()
$ws.Range("E38").Value = '  +2.86%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '8.66'
$c.ClearFormats()
$ws.Range("E39").Value = '  -0.11%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.79'
$c.ClearFormats()
$ws.Range("E40").Value = '  +10.02%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '402.88'
$c.ClearFormats()
$ws.Range("E41").Value = '  +7.77%  '
$ws.Range("E42").Value = '  +3.49%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.784.13'
$c.ClearFormats()
$ws.Range("E43").Value = '  +4.86%  '
$ws.Range("E44").Value = '  +0.74%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.254'
$c.ClearFormats()
$ws.Range("E45").Value = '  +7.37%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '123.25'
$c.ClearFormats()
$ws.Range("E47").Value = '  +1.02%  '
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("E49").Value = '  +3.42%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '33.68'
$c.ClearFormats()
$ws.Range("E50").Value = '  +21.22%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '23.53'
$c.ClearFormats()
$ws.Range("E51").Value = '  +2.12%  '
